$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 5) that mirrors the existing rows 2-4,
# with the new values from the commit.
$ws.Range("A5").Value = "N"
$ws.Range("B5").Value = "ConversionRate"
$ws.Range("C5").Value = "CurrencyConvertorSoap"
$ws.Range("D5").Value = "http://www.webservicex.com/currencyconvertor.asmx"
$ws.Range("E5").Value = "SIT"
$ws.Range("F5").Value = "Test4"
$ws.Range("G5").Value = "USD"
$ws.Range("H5").Value = "DKK"

# Copy the formatting (border style) from row 4 onto the new row 5 so the
# new cells share the existing bordered style instead of the default one.
$ws.Range("A4:H4").Copy()
$ws.Range("A5:H5").PasteSpecial(-4122)

# Add the hyperlink for D5 pointing at the same endpoint URL used by the
# other rows (this registers a new rId4 relationship).
$ws.Hyperlinks.Add($ws.Range("D5"), "http://www.webservicex.com/currencyconvertor.asmx")

# Hyperlinks.Add re-styles the cell with the built-in "Hyperlink" look
# (blue/underline); restore the same plain bordered style the other
# hyperlinked cells in column D use.
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)

# Match the saved selection/active cell recorded in the workbook.
$ws.Range("H9").Select()
